# "Add files via upload" -- add CZ error and CZ cycle error to 10q exp parameters
#
# The workbook originally had one sheet named "CZ" holding the per-pair
# error numbers that are actually the CZ *cycle* errors. This edit:
#   1. Renames that sheet "CZ" -> "CZ_cycle" (its data/labels are untouched;
#      only the column header text becomes "cz_errors_cycle").
#   2. Inserts a brand-new sheet named "CZ" right after "CZ_cycle",
#      containing the newly measured per-pair CZ gate error numbers under
#      the header "cz_errors", using the same qubit-pair row labels as the
#      CZ_cycle sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the existing "CZ" sheet to "CZ_cycle" ------------------------
$czCycleSheet = $wb.Worksheets.Item("CZ")
$czCycleSheet.Name = "CZ_cycle"
$czCycleSheet.Range("B1").Value = "cz_errors_cycle"

# --- 2. Insert the new "CZ" sheet right after "CZ_cycle" --------------------
$newCzSheet = $wb.Worksheets.Add($null, $czCycleSheet)
$newCzSheet.Name = "CZ"

# Header cell - copy first so it inherits the bold/border/center style used
# by every other sheet's header row, then overwrite the text.
$czCycleSheet.Range("B1").Copy($newCzSheet.Range("B1"))
$newCzSheet.Range("B1").Value = "cz_errors"

# Row labels (same qubit pairs as CZ_cycle, same styling) copied in one shot.
$czCycleSheet.Range("A2:B10").Copy($newCzSheet.Range("A2"))

# New CZ gate error values measured for each qubit pair.
$values = @(
    0.002387399578551386,
    0.002228112826128892,
    0.002137539206894301,
    0.005197383485476048,
    0.005299221208069804,
    0.002660335286164495,
    0.003823982139691484,
    0.002966396423230599,
    0.003615822736568575
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $newCzSheet.Cells.Item($row, 2).Value = $values[$i]
}

# Restore the originally-active sheet/selection.
$wb.Worksheets.Item("readout_errors").Activate()
